$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Extend header row (row 1) with two new columns P1, Q1 ---
# Copy formatting from O1 (existing last header cell: bold font, thin box border,
# centered horizontal / top vertical alignment) onto the two new header cells,
# then set their values.
$ws.Range("O1").Copy() | Out-Null
$ws.Range("P1:Q1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15

# --- Update data rows 2 through 25 ---
for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 9).Value = 2   # column I -> 2
    $ws.Cells.Item($r, 11).Value = 1  # column K -> 1
    $ws.Cells.Item($r, 13).Value = 2  # column M -> 2
    $ws.Cells.Item($r, 15).Value = 1  # column O -> 1

    $ws.Cells.Item($r, 16).Value = 2  # column P -> 2 (new)
    $ws.Cells.Item($r, 17).Value = 2  # column Q -> 2 (new)
}
